$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J header ("Misc. Data")
$ws.Range("J1").Value = 'Misc. Data'

# Row 2
$ws.Range("E2").Value = '[Tran Duc Anh%Ly%NULL%3,                         Didier%Zanini%NULL%1,                         Vincent%Laforge%NULL%1,                         Sylvie%Arlotto%NULL%1,                         Stephanie%Gentile%NULL%1,                         Helene%Mendizabal%NULL%1,                         Michael%Finaud%NULL%1,                         David%Morel%NULL%1,                         Olivier%Quenette%NULL%1,                         Priscilla%Malfuson-Clot-Faybesse%NULL%1,                         Alain%Midejean%NULL%1,                         Phuc%Le-Dinh%NULL%1,                         Gérard%Daher%NULL%1,                         Berengere%Labarriere%NULL%1,                         Anne-Marie%Morel-Roux%NULL%1,                         Alain%Coquet%NULL%1,                         Patrick%Augier%NULL%1,                         Philippe%Parola%NULL%2,                         Eric%Chabriere%NULL%1,                         Didier%Raoult%NULL%2,                         Philippe%Gautret%NULL%1]'
$ws.Range("J2").Value = ''

# Row 3
$ws.Range("E3").Value = '[Matthieu%Mahévas%NULL%0,                         Viet-Thi%Tran%NULL%0,                         Viet-Thi%Tran%NULL%0,                         Mathilde%Roumier%NULL%0,                         Amélie%Chabrol%NULL%0,                         Romain%Paule%NULL%0,                         Constance%Guillaud%NULL%0,                         Elena%Fois%NULL%0,                         Raphael%Lepeule%NULL%0,                         Tali-Anne%Szwebel%NULL%0,                         François-Xavier%Lescure%NULL%0,                         Frédéric%Schlemmer%NULL%0,                         Marie%Matignon%NULL%0,                         Mehdi%Khellaf%NULL%0,                         Etienne%Crickx%NULL%0,                         Benjamin%Terrier%NULL%0,                         Caroline%Morbieu%NULL%0,                         Paul%Legendre%NULL%0,                         Julien%Dang%NULL%0,                         Yoland%Schoindre%NULL%0,                         Jean-Michel%Pawlotsky%NULL%0,                         Marc%Michel%NULL%0,                         Elodie%Perrodeau%NULL%0,                         Nicolas%Carlier%NULL%0,                         Nicolas%Roche%NULL%0,                         Victoire%de Lastours%NULL%0,                         Clément%Ourghanlian%NULL%0,                         Solen%Kerneis%NULL%0,                         Philippe%Ménager%NULL%0,                         Luc%Mouthon%NULL%0,                         Etienne%Audureau%NULL%0,                         Philippe%Ravaud%NULL%0,                         Bertrand%Godeau%NULL%0,                         Sébastien%Gallien%NULL%0,                         Nathalie%Costedoat-Chalumeau%NULL%0]'
$ws.Range("J3").Value = ''

# Row 4
$ws.Range("E4").Value = '[Ahmed%Hasan%coreGivesNoEmail%1,                        Amir%Lerman%coreGivesNoEmail%1,                        Charanjit%Rihal%coreGivesNoEmail%1,                        Charles%Cagin%coreGivesNoEmail%1,                        Derek%So%coreGivesNoEmail%1,                        Erin%Iturriaga%coreGivesNoEmail%1,                        Ivan%Chavez%coreGivesNoEmail%1,                        J.%Dawn Abbott%coreGivesNoEmail%1,                        Jang-Ho%Bae%coreGivesNoEmail%1,                        Jean-Francois%Tanguay%coreGivesNoEmail%1,                        Kent%Bailey%coreGivesNoEmail%1,                        Liewei%Wang%coreGivesNoEmail%1,                        Linnea%Baudhuin%coreGivesNoEmail%1,                        Malcolm%Bell%coreGivesNoEmail%1,                        Mandeep%Sidhu%coreGivesNoEmail%1,                        Michael%E. Farkouh%coreGivesNoEmail%1,                        Myung%Ho Jeong%coreGivesNoEmail%1,                        Nancy%Geller%coreGivesNoEmail%1,                        Naveen%L. Pereira%coreGivesNoEmail%1,                        Paul%Gordon%coreGivesNoEmail%1,                        Richard%Weinshilboum%coreGivesNoEmail%1,                        Robert%Welsh%coreGivesNoEmail%1,                        Ryan%Lennon%coreGivesNoEmail%1,                        Shaun%G. Goodman%coreGivesNoEmail%1,                        Verghese%Mathew%coreGivesNoEmail%1,                        Yi-Ping%Fu%coreGivesNoEmail%1,                        Yves%Rosenberg%coreGivesNoEmail%1]'
$ws.Range("J4").Value = ''

# Row 5 (E5 unchanged)
$ws.Range("J5").Value = ''

# Row 6
$ws.Range("E6").Value = '[Francisco Javier%Membrillo de Novales%xref no email%1,     Germ\u00e1n%Ram\u00edrez-Olivencia%xref no email%1,     Miriam%Est\u00e9banez%xref no email%1,     Bego\u00f1a%de Dios%xref no email%1,     Mar\u00eda Dolores%Herrero%xref no email%1,     Tatiana%Mata%xref no email%1,     Alberto M.%Borobia%xref no email%1,     Carlos%Guti\u00e9rrez%xref no email%1,     Mar\u00eda%Sim\u00f3n%xref no email%1,     Ana%Ochoa%xref no email%1,     Yolanda%Mart\u00ednez%xref no email%1,     Alejandro%Aguirre%xref no email%1,     Francisco de As\u00eds%Alc\u00e1ntara%xref no email%1,     Pablo%Fern\u00e1ndez-Gonz\u00e1lez%xref no email%1,     Elena%L\u00f3pez%xref no email%1,     Sergio%Campos%xref no email%1,     Mar\u00eda%Navarro%xref no email%1,     Luc\u00eda Elena%Ballester%xref no email%1]'
$ws.Range("I6").Value = ''
$ws.Range("J6").Value = 'PUBLISHER: MDPI AG'

# Row 7
$ws.Range("E7").Value = '[Joseph%Magagnoli%NULL%1,                         Siddharth%Narendran%NULL%1,                         Felipe%Pereira%NULL%1,                         Tammy H.%Cummings%NULL%1,                         James W.%Hardin%NULL%1,                         S. Scott%Sutton%NULL%1,                         Jayakrishna%Ambati%NULL%1]'
$ws.Range("J7").Value = ''

# Row 8
$ws.Range("E8").Value = '[ S.%Singh%null%1,                          A.% Khan%null%1,                          M.% Chowdhry%null%1,                          A. % Chatterjee%null%1,       S.%Singh%null%1,       A.% Khan%null%1,       M.% Chowdhry%null%1,       A. % Chatterjee%null%1]'
$ws.Range("J8").Value = ''

